$d = $word.ActiveDocument

# Paragraph 2
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="065A3508" w14:textId="6D8486BC" w:rsidR="00F94FF0" w:rsidRDefault="0078168E"><w:r><w:t xml:space="preserve">You will be creating </w:t></w:r><w:r w:rsidR="009A3539"><w:t>two</w:t></w:r><w:r w:rsidR="00F94FF0"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>w</w:t></w:r><w:r><w:t>eb pages described in Tutorial 2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> of our textbook, </w:t></w:r><w:r w:rsidRPr="00CE6DC1"><w:rPr><w:i/></w:rPr><w:t>HTML, CSS and Dynamic HTML</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00F00B3A"><w:t>For best results, read</w:t></w:r><w:r w:rsidR="00F94FF0"><w:t xml:space="preserve"> the</w:t></w:r><w:r w:rsidR="00F00B3A"><w:t xml:space="preserve"> explanations in the tutorial (</w:t></w:r><w:r w:rsidR="00267155"><w:t xml:space="preserve">but </w:t></w:r><w:r w:rsidR="00F00B3A"><w:t>not the steps) before writing HTML code</w:t></w:r><w:r w:rsidR="00F94FF0"><w:t xml:space="preserve"> for your web pages</w:t></w:r><w:r w:rsidR="00F00B3A"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="00267155"><w:t xml:space="preserve">You might want to read a few pages and then try writing some code, read a few pages, </w:t></w:r><w:r w:rsidR="00F94FF0"><w:t xml:space="preserve">then </w:t></w:r><w:r w:rsidR="00267155"><w:t xml:space="preserve">write some more code. </w:t></w:r></w:p>'
$d.Paragraphs(2).Range.InsertXML($xml2)

# Paragraph 6
$xml6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="496060AE" w14:textId="20C061F9" w:rsidR="00CE6DC1" w:rsidRDefault="00F94FF0" w:rsidP="0078168E"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="00F94FF0"><w:rPr><w:b/></w:rPr><w:t>Tutorial</w:t></w:r><w:r w:rsidRPr="00F94FF0"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> 2</w:t></w:r><w:r w:rsidRPr="00F94FF0"><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0078168E"><w:t>C</w:t></w:r><w:r w:rsidR="00CE6DC1"><w:t xml:space="preserve">reate the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00466CF4"><w:rPr><w:i/></w:rPr><w:t>CAMshots</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0078168E" w:rsidRPr="00F94FF0"><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0078168E"><w:t>web page</w:t></w:r><w:r w:rsidR="00CE6DC1"><w:t xml:space="preserve"> described in t</w:t></w:r><w:r w:rsidR="0078168E"><w:t xml:space="preserve">he tutorial. </w:t></w:r><w:r><w:br/></w:r><w:r w:rsidR="0078168E"><w:br/></w:r><w:r w:rsidR="00F00B3A" w:rsidRPr="00F94FF0"><w:rPr><w:b/></w:rPr><w:t>Best Practice</w:t></w:r><w:r><w:t>: M</w:t></w:r><w:r w:rsidR="0078168E"><w:t>ake the web page</w:t></w:r><w:r w:rsidR="00466CF4"><w:t xml:space="preserve"> by just looking at the pictures of the completed pages </w:t></w:r><w:r w:rsidR="00F00B3A"><w:t>and then trying</w:t></w:r><w:r w:rsidR="0078168E"><w:t xml:space="preserve"> to write HTML that make</w:t></w:r><w:r w:rsidR="00F00B3A"><w:t>s</w:t></w:r><w:r w:rsidR="0078168E"><w:t xml:space="preserve"> a web page </w:t></w:r><w:r w:rsidR="00466CF4"><w:t>look</w:t></w:r><w:r w:rsidR="0078168E"><w:t xml:space="preserve"> just like the picture. If you run into difficulty, then look at the step by step instructions.</w:t></w:r><w:r w:rsidR="0096682F"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00466CF4"><w:t xml:space="preserve">In “real life” </w:t></w:r><w:r><w:t>you wouldn’t have step by step instructions.</w:t></w:r><w:r><w:br/></w:r></w:p>'
$d.Paragraphs(6).Range.InsertXML($xml6)

# Paragraph 7
$xml7 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1169CEC7" w14:textId="34C9FC0D" w:rsidR="0096682F" w:rsidRDefault="00F94FF0" w:rsidP="0096682F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r w:rsidRPr="00F94FF0"><w:rPr><w:b/></w:rPr><w:t>Review:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00466CF4"><w:t>Add the Glossary</w:t></w:r><w:r><w:t xml:space="preserve"> web page</w:t></w:r><w:r w:rsidR="00466CF4"><w:t xml:space="preserve"> to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00466CF4"><w:t>CAMshots</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00466CF4"><w:t xml:space="preserve"> site</w:t></w:r><w:r><w:t xml:space="preserve">. Again, it’s best to just look at the picture </w:t></w:r><w:r w:rsidR="00466CF4"><w:t>of the completed page</w:t></w:r><w:r><w:t xml:space="preserve"> and write HTML code to make a web page that looks just like the picture. </w:t></w:r></w:p>'
$d.Paragraphs(7).Range.InsertXML($xml7)

# Paragraph 11
$xml11 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4860E211" w14:textId="73477901" w:rsidR="00882812" w:rsidRDefault="00710B65" w:rsidP="0096682F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Zip the tutorial.02</w:t></w:r><w:r w:rsidR="0096682F"><w:t xml:space="preserve"> folder along with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="0096682F"><w:t>it’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="0096682F"><w:t xml:space="preserve"> tutorial, review, and case subfolders.</w:t></w:r></w:p>'
$d.Paragraphs(11).Range.InsertXML($xml11)

# Paragraph 12
$xml12 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7F00BB4E" w14:textId="71F1D08E" w:rsidR="0096682F" w:rsidRDefault="00710B65" w:rsidP="0096682F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Upload your tutorial.02</w:t></w:r><w:r w:rsidR="00355E86"><w:t xml:space="preserve">.zip file using the </w:t></w:r><w:r w:rsidR="0096682F" w:rsidRPr="00413A30"><w:rPr><w:i/></w:rPr><w:t>Production Version</w:t></w:r><w:r w:rsidR="0096682F"><w:t xml:space="preserve"> assignment</w:t></w:r><w:r w:rsidR="00355E86"><w:t xml:space="preserve"> link on Moodle</w:t></w:r><w:r w:rsidR="0096682F"><w:t>:</w:t></w:r></w:p>'
$d.Paragraphs(12).Range.InsertXML($xml12)

# Header paragraph
$hdrXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="339DD4E7" w14:textId="004C2D65" w:rsidR="00A02234" w:rsidRPr="00A02234" w:rsidRDefault="00A02234" w:rsidP="00A02234"><w:pPr><w:pStyle w:val="Header"/><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Instructions for </w:t></w:r><w:r w:rsidR="00466CF4"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Lab 3</w:t></w:r><w:r w:rsidR="00F94FF0"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00466CF4"><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>CAMshots</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:br/></w:r><w:r w:rsidRPr="00A02234"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>CIS 195 Web Authoring 1</w:t></w:r></w:p>'
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Paragraphs(1).Range.InsertXML($hdrXml)

Write-Output "done"